$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing A:D data shifts to B:E
$ws.Range("A1").EntireColumn.Insert()

# New header cell for the inserted column
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats - reuse the header style

# Fill new column A with the 0-based segment index (same style as the label column B)
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($row, 1).Value = $i
}

$excel.CutCopyMode = 0
